# Generate Report for Handback
# Update the "last generated" timestamps recorded on the handback status
# report. These values are stored as plain text (not Excel dates), so we
# must assign them as strings to avoid any numeric/date auto-conversion.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the
# 47c1bbe3-...-8dfe9b704f15.md row (row 3).
# This value is shared with the de-de sheet's "Correspond Handoff Datetime"
# for the same source file, so both cells are updated together.
$wsOverview.Range("G3").Value = "2016-09-05 11:05:55"
$wsDeDe.Range("H3").Value = "2016-09-05 11:05:55"

# zh-cn sheet, row 3 (47c1bbe3-...-8dfe9b704f15.md):
#   Correspond Handoff Datetime
$wsZhCn.Range("H3").Value = "2016-09-05 11:05:49"
#   Correspond Handback DateTime
$wsZhCn.Range("K3").Value = "2016-09-05 11:06:29"

# de-de sheet, row 3 (47c1bbe3-...-8dfe9b704f15.md):
#   Correspond Handback DateTime
$wsDeDe.Range("K3").Value = "2016-09-05 11:06:37"
